$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: O31 flips from "√" to "×"
$ws.Range("O31").Value = "×"

# Row 32: finish filling out 2019-07-31 -> becomes 2019-08-05 entry
$ws.Range("A32").Value = 20190805
$ws.Range("B32:C32").Value = "√"
$ws.Range("D32").Value = "×"
$ws.Range("E32:O32").Value = "√"

# Row 33: 2019-08-06
$ws.Range("A33").Value = 20190806
$ws.Range("B33:C33").Value = "√"
$ws.Range("D33").Value = "×"
$ws.Range("E33:O33").Value = "√"

# Row 34: 2019-08-07
$ws.Range("A34").Value = 20190807
$ws.Range("B34:G34").Value = "√"
$ws.Range("H34").Value = "×"
$ws.Range("I34:N34").Value = "√"
$ws.Range("O34").Value = "×"

# Row 35: 2019-08-08 (only started, A-C filled so far)
$ws.Range("A35").Value = 20190808
$ws.Range("B35:C35").Value = "√"

$ws.Range("C35").Select()
